$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.045.07'
$ws.Range("E2").Value = '  +1.85%  '

$ws.Range("D3").Value = '3.833.21'
$ws.Range("E3").Value = '  +0.57%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.45%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '625.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +3.74%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '165.39'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.51%  '

$ws.Range("D7").Value = '3.826.08'
$ws.Range("E7").Value = '  +0.41%  '

$ws.Range("E8").Value = '  +0.00%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.519'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.161'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.96%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.454'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +0.13%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.68'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +3.51%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000249'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '35.78'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.74%  '

$ws.Range("D15").Value = '4.461.37'
$ws.Range("E15").Value = '  +0.32%  '

$ws.Range("D16").Value = '3.809.77'
$ws.Range("E16").Value = '  +0.27%  '

$ws.Range("D17").Value = '68.947.53'
$ws.Range("E17").Value = '  +1.68%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.18'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -1.39%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.14'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.64%  '

$ws.Range("E20").Value = '  -0.24%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '468.77'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.03%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.75'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.77%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.706'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.26%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000151'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +2.43%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '84.07'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '12.07'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.15'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +1.55%  '

$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.51%  '

$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("D30").Value = '3.968.88'
$ws.Range("E30").Value = '  +0.20%  '

$ws.Range("B31").Value = 'ImmutableX'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.25'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +1.80%  '

$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.66'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.22%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.29'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.78%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.26'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -0.55%  '

$ws.Range("B35").Value = 'Aptos'
$ws.Range("C35").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.10'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +0.27%  '

$ws.Range("B36").Value = 'Binance-PegBSC-USD'
$ws.Range("C36").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.998'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.29%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.102'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.59%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.148'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +6.97%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.92'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +1.70%  '

$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.24'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("B41").Value = 'Mantle'
$ws.Range("C41").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.982'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -1.74%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.998'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '157.31'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  +3.67%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.300'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.51%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '43.05'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -4.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.41'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +2.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '46.74'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.06%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.41'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.46%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.89'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.44%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '381.27'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.35%  '
